# "new format for excel"
# The Public/Private/Save/Cache/Ref/Upload/Desc flag columns (B:F) for rows
# 3-9 were stored as Boolean (TRUE/FALSE) cells. Re-write them as plain
# numbers (0) instead, which is the new format requested by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3:F9").Value = 0
